# econ: anchor machines for wheeled harvester and forwarders
#
# Inserts two new parameter columns on the "parameterization" sheet:
#   - addOnWinchCableLength (before the old ctlHaulHours column, i.e. BF)
#   - anchorSMh (before the old grappleYardingConstant column, i.e. CD,
#     which is CE once the first insertion has shifted things right)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameterization")

# --- Insert both new columns first (shifts everything else right) -------
$ws.Range("BF1").EntireColumn.Insert()
$ws.Range("CE1").EntireColumn.Insert()

# --- Populate them. "anchorSMh" must be registered as a shared string
# before "addOnWinchCableLength" to match the original authoring order.
$ws.Range("CE1").Value = "anchorSMh"
$ws.Range("CE2").Value = 71.5

$ws.Range("BF1").Value = "addOnWinchCableLength"
$ws.Range("BF2").Value = 350
$ws.Range("BF2").NumberFormat = "0"

# --- Restore the view state (frozen pane / selection) --------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 55
$ws.Range("BF1:BF2").Select()
